# Adds a new "Sim1 N=4000" worksheet (a second simulation run, this time with
# N=4000 total observations) positioned between "Sim1" and "Sheet1", mirroring
# the layout/formatting of "Sim1" but with the smaller-N results and becomes the
# active tab, matching commit "Simulation 1 (10 studies): Results".

$wb = $excel.ActiveWorkbook
$sim1 = $wb.Worksheets.Item("Sim1")

# Duplicate "Sim1" right after itself -- this clones all values/number formats/
# column widths/row heights, which we then touch up below.
$sim1.Copy($null, $sim1)
$ws = $wb.Worksheets.Item(2)
$ws.Name = "Sim1 N=4000"

# --- Update the Study.1..Study.10 results for the N=4000 simulation run ---
$ws.Range("C2").Value = 100
$ws.Range("D2").Value = 100
$ws.Range("E2").Value = 100
$ws.Range("F2").Value = 100
$ws.Range("G2").Value = 100
$ws.Range("H2").Value = 100
$ws.Range("I2").Value = 100
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 100
$ws.Range("L2").Value = 100
$ws.Range("C3").Value = 50
$ws.Range("D3").Value = 50
$ws.Range("E3").Value = 50
$ws.Range("F3").Value = 50
$ws.Range("G3").Value = 50
$ws.Range("H3").Value = 150
$ws.Range("I3").Value = 150
$ws.Range("J3").Value = 150
$ws.Range("K3").Value = 150
$ws.Range("L3").Value = 150
$ws.Range("C4").Value = 25
$ws.Range("D4").Value = 25
$ws.Range("E4").Value = 25
$ws.Range("F4").Value = 25
$ws.Range("G4").Value = 25
$ws.Range("H4").Value = 175
$ws.Range("I4").Value = 175
$ws.Range("J4").Value = 175
$ws.Range("K4").Value = 175
$ws.Range("L4").Value = 175
$ws.Range("F5").Value = 100
$ws.Range("G5").Value = 100
$ws.Range("H5").Value = 100
$ws.Range("I5").Value = 100
$ws.Range("J5").Value = 175
$ws.Range("K5").Value = 175
$ws.Range("L5").Value = 175
$ws.Range("F6").Value = 45
$ws.Range("G6").Value = 45
$ws.Range("H6").Value = 45
$ws.Range("I6").Value = 45
$ws.Range("J6").Value = 245
$ws.Range("K6").Value = 250
$ws.Range("L6").Value = 250
$ws.Range("D7").Value = 50
$ws.Range("E7").Value = 50
$ws.Range("F7").Value = 50
$ws.Range("G7").Value = 50
$ws.Range("H7").Value = 100
$ws.Range("I7").Value = 100
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 100
$ws.Range("L7").Value = 375
$ws.Range("L8").Value = 375
$ws.Range("C9").Value = 25
$ws.Range("D9").Value = 25
$ws.Range("E9").Value = 25
$ws.Range("F9").Value = 25
$ws.Range("G9").Value = 25
$ws.Range("H9").Value = 25
$ws.Range("I9").Value = 25
$ws.Range("J9").Value = 25
$ws.Range("K9").Value = 25
$ws.Range("L9").Value = 775
$ws.Range("E10").Value = 25
$ws.Range("F10").Value = 25
$ws.Range("G10").Value = 50
$ws.Range("H10").Value = 50
$ws.Range("I10").Value = 100
$ws.Range("J10").Value = 100
$ws.Range("K10").Value = 200
$ws.Range("L10").Value = 400

# --- Re-colour the handful of cells whose heat-map shading changed tone ---
$ws.Range("C2").Copy()
$ws.Range("D2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C2").Copy()
$ws.Range("E2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C2").Copy()
$ws.Range("F2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C2").Copy()
$ws.Range("G2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C2").Copy()
$ws.Range("I2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C2").Copy()
$ws.Range("J2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C2").Copy()
$ws.Range("K2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("H3").Copy()
$ws.Range("J3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("H3").Copy()
$ws.Range("K3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C4").Copy()
$ws.Range("D4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C4").Copy()
$ws.Range("E4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C4").Copy()
$ws.Range("F4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C4").Copy()
$ws.Range("G4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("H3").Copy()
$ws.Range("I4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("H3").Copy()
$ws.Range("J4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("H3").Copy()
$ws.Range("K4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C2").Copy()
$ws.Range("I5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("J5").Copy()
$ws.Range("L5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C4").Copy()
$ws.Range("I6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("J6").Copy()
$ws.Range("L6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("L7").Copy()
$ws.Range("L8").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C5").Copy()
$ws.Range("E10").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C5").Copy()
$ws.Range("F10").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Restore the author's last selection and make the new sheet the active tab
$ws.Range("F14").Select() | Out-Null
$ws.Activate() | Out-Null
